# Apply updated Samsung Phones Data rows (price list refresh + two new rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

Set-TextCell 2 1 'Samsung Galaxy M32 5G (Sky Blue, 6GB RAM, 128GB Storage)'
Set-TextCell 2 2 '16,999'
Set-TextCell 3 1 'Samsung Galaxy Note 20 (Mystic Green, 8GB RAM, 256GB Storage) with No Cost EMI/Additional Exchange Offers'
Set-TextCell 3 2 '44,999'
Set-TextCell 4 1 'Samsung Galaxy M32 5G (Sky Blue, 8GB RAM, 128GB Storage)'
Set-TextCell 4 2 '18,999'
Set-TextCell 5 1 'Samsung Galaxy M31 (Ocean Blue, 8GB RAM, 128GB Storage) 6 Months Free Screen Replacement for Prime'
Set-TextCell 5 2 '15,999'
Set-TextCell 6 1 'Samsung Galaxy M32 5G (Slate Black, 8GB RAM, 128GB Storage)'
Set-TextCell 6 2 '18,999'
Set-TextCell 7 1 'Samsung Galaxy M32 (Light Blue, 6GB RAM, 128GB Storage) 6 Months Free Screen Replacement for Prime'
Set-TextCell 7 2 '14,999'
Set-TextCell 8 1 'Samsung Galaxy M32 (Black, 4GB RAM, 64GB Storage) 6 Months Free Screen Replacement for Prime'
Set-TextCell 8 2 '12,999'
Set-TextCell 9 1 'Samsung Galaxy M12 (Blue,4GB RAM, 64GB Storage) 6000 mAh with 8nm Processor | True 48 MP Quad Camera | 90Hz Refresh Rate'
Set-TextCell 9 2 '9,499'
Set-TextCell 10 1 'Samsung Galaxy M21 2021 Edition (Arctic Blue, 4GB RAM, 64GB Storage) | FHD+ sAMOLED | 6 Months Free Screen Replacement for Prime (SM-M215GLBDINS)'
Set-TextCell 10 2 '11,999'
Set-TextCell 11 1 'Samsung Galaxy M51 (Celestial Black, 6GB RAM, 128GB Storage) 6 Months Free Screen Replacement for Prime'
Set-TextCell 11 2 '19,999'
Set-TextCell 12 1 'Samsung Galaxy M32 5G (Sky Blue, 6GB RAM, 128GB Storage)'
Set-TextCell 12 2 '16,999'
Set-TextCell 13 1 'Samsung Galaxy M12 (Black,6GB RAM, 128GB Storage) 6 Months Free Screen Replacement for Prime'
Set-TextCell 13 2 '11,499'
Set-TextCell 14 1 'Samsung Galaxy M12 (Black,4GB RAM, 64GB Storage) 6000 mAh with 8nm Processor | True 48 MP Quad Camera | 90Hz Refresh Rate'
Set-TextCell 14 2 '9,499'
Set-TextCell 15 1 'Samsung Galaxy M32 (Light Blue, 4GB RAM, 64GB Storage) 6 Months Free Screen Replacement for Prime'
Set-TextCell 15 2 '12,999'
Set-TextCell 16 1 'Samsung Galaxy M32 (Black, 6GB RAM, 128GB Storage) 6 Months Free Screen Replacement for Prime'
Set-TextCell 16 2 '14,999'
Set-TextCell 17 1 'Samsung Galaxy M52 5G (Blazing Black, 6GB RAM, 128GB Storage) Latest Snapdragon 778G 5G | sAMOLED 120Hz Display'
Set-TextCell 17 2 '25,999'
Set-TextCell 18 1 'Samsung Galaxy M12 (Blue,6GB RAM, 128GB Storage) 6 Months Free Screen Replacement for Prime'
Set-TextCell 18 2 '11,499'
Set-TextCell 19 1 'Samsung Galaxy M12 (White,6GB RAM, 128GB Storage) 6 Months Free Screen Replacement for Prime'
Set-TextCell 19 2 '11,499'
Set-TextCell 20 1 'Samsung Galaxy S20 FE 5G (Cloud Navy, 8GB RAM, 128GB Storage)'
Set-TextCell 20 2 '36,990'
Set-TextCell 21 1 'Samsung Galaxy M21 2021 Edition (Charcoal Black, 4GB RAM, 64GB Storage) | FHD+ sAMOLED | 6 Months Free Screen Replacement for Prime (SM-M215GZKDINS)'
Set-TextCell 21 2 '11,999'
Set-TextCell 22 1 'Samsung Galaxy M21 2021 Edition (Arctic Blue, 6GB RAM, 128GB Storage) | FHD+ sAMOLED | 6 Months Free Screen Replacement for Prime (SM-M215GLBHINS)'
Set-TextCell 22 2 '13,999'
Set-TextCell 23 1 'Samsung Galaxy M31 (Space Black, 8GB RAM, 128GB Storage) 6 Months Free Screen Replacement for Prime'
Set-TextCell 23 2 '15,999'
Set-TextCell 24 1 'Samsung Galaxy M52 5G (Blazing Black, 8GB RAM, 128GB Storage) Latest Snapdragon 778G 5G | sAMOLED 120Hz Display'
Set-TextCell 24 2 '27,999'
Set-TextCell 25 1 'Samsung Galaxy M51 (Electric Blue, 6GB RAM, 128GB Storage) 6 Months Free Screen Replacement for Prime'
Set-TextCell 25 2 '19,999'
Set-TextCell 26 1 'Samsung Galaxy A22 (Black, 6GB RAM, 128GB Storage) with No Cost EMI/Additional Exchange Offers'
Set-TextCell 26 2 '18,499'
Set-TextCell 27 1 'Samsung Galaxy M52 5G (ICY Blue, 6GB RAM, 128GB Storage) Latest Snapdragon 778G 5G | sAMOLED 120Hz Display'
Set-TextCell 27 2 '25,999'
